$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.151.61"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "3.148.97"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'569.21"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'162.49"
$ws.Range("E6").Value = "  -3.96%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.574"
$ws.Range("E8").Value = "  -5.68%  "
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("D10").Value = "'6.58"
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").Value = "'0.380"
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("D12").Value = "3.702.91"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").Value = "64.229.85"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "'24.99"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").Value = "3.147.13"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("E17").Value = "  -2.83%  "
$ws.Range("D18").Value = "'402.57"
$ws.Range("E18").Value = "  -4.23%  "
$ws.Range("D19").Value = "'12.60"
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").Value = "'5.21"
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("D21").Value = "'7.08"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  +3.48%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'68.40"
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("D25").Value = "'0.481"
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("E26").Value = "  -5.21%  "
$ws.Range("E27").Value = "  -5.22%  "
$ws.Range("D28").Value = "'8.78"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("D31").Value = "'21.08"
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("D32").Value = "'6.25"
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("D33").Value = "'4.80"
$ws.Range("E33").Value = "  -4.77%  "
$ws.Range("D34").Value = "'156.08"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("E36").Value = "  -3.49%  "
$ws.Range("D37").Value = "2.660.83"
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("D39").Value = "'23.58"
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("D40").Value = "'4.06"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").Value = "'0.691"
$ws.Range("E41").Value = "  -2.68%  "
$ws.Range("D42").Value = "'0.0613"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").Value = "'5.38"
$ws.Range("E43").Value = "  -5.76%  "
$ws.Range("D44").Value = "'0.0254"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").Value = "'287.16"
$ws.Range("E45").Value = "  -3.33%  "
$ws.Range("D46").Value = "'21.09"
$ws.Range("E46").Value = "  -3.83%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").Value = "'10.50"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").Value = "'1.87"
$ws.Range("E50").Value = "  -8.69%  "
$ws.Range("D51").Value = "'5.65"
$ws.Range("E51").Value = "  -2.28%  "

# Reset style on cells forced to text via apostrophe prefix so no new
# number-format style entries are introduced (cells had no explicit style).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
